# Update for release to deploy 0.1.1
$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from NMDP Disease Cod" tab to "Include #0"
$includeSheet = $wb.Worksheets.Item("Include from NMDP Disease Cod")
$includeSheet.Name = "Include #0"

# 2. Update the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")

# Bump version
$meta.Range("B3").Value = "0.1.1"

# Bump the date stamp
$meta.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new "Jurisdiction" row (with an empty value) right after "Contact" (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row. All of these data
# rows already share identical formatting, so extend that formatting onto the new last
# row (15), then shift the values down from the bottom up rather than using
# Rows.Insert() (which would otherwise synthesize new/duplicate style records).
$meta.Range("A14:B14").Copy()
$meta.Range("A15:B15").PasteSpecial(-4122)

$meta.Range("A15").Value = $meta.Range("A14").Value2
$meta.Range("B15").Value = $meta.Range("B14").Value2

$meta.Range("A14").Value = $meta.Range("A13").Value2
$meta.Range("B14").Value = $meta.Range("B13").Value2

$meta.Range("A13").Value = $meta.Range("A12").Value2
$meta.Range("B13").Value = $meta.Range("B12").Value2

$meta.Range("A12").Value = $meta.Range("A11").Value2
$meta.Range("B12").Value = $meta.Range("B11").Value2

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""
